# Weekly update: insert a new week's Apio (Primera/Segunda) price entries
# above the existing row 120, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 120-121; everything currently at row 120
# onward (through the old row 192) shifts down to rows 122-194.
$ws.Rows("120:121").Insert()

# New row 120: Apio, Americana (o), Primera
$ws.Cells.Item(120, 1).Value = 11
$ws.Cells.Item(120, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(120, 3).Value = "Bíobío"
$ws.Cells.Item(120, 4).Value = 44582
$ws.Cells.Item(120, 5).Value = 8
$ws.Cells.Item(120, 6).Value = 100112017
$ws.Cells.Item(120, 7).Value = "Apio"
$ws.Cells.Item(120, 8).Value = "Americana (o)"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 100
$ws.Cells.Item(120, 11).Value = 6500
$ws.Cells.Item(120, 12).Value = 7000
$ws.Cells.Item(120, 13).Value = 6750
$ws.Cells.Item(120, 14).Value = "$/docena de matas"
$ws.Cells.Item(120, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(120, 16).Value = 1125
$ws.Cells.Item(120, 17).Value = 6
$ws.Cells.Item(120, 18).Value = "Hortaliza"

# New row 121: Apio, Americana (o), Segunda
$ws.Cells.Item(121, 1).Value = 11
$ws.Cells.Item(121, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(121, 3).Value = "Bíobío"
$ws.Cells.Item(121, 4).Value = 44582
$ws.Cells.Item(121, 5).Value = 8
$ws.Cells.Item(121, 6).Value = 100112017
$ws.Cells.Item(121, 7).Value = "Apio"
$ws.Cells.Item(121, 8).Value = "Americana (o)"
$ws.Cells.Item(121, 9).Value = "Segunda"
$ws.Cells.Item(121, 10).Value = 50
$ws.Cells.Item(121, 11).Value = 6000
$ws.Cells.Item(121, 12).Value = 6000
$ws.Cells.Item(121, 13).Value = 6000
$ws.Cells.Item(121, 14).Value = "$/docena de matas"
$ws.Cells.Item(121, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(121, 16).Value = 1000
$ws.Cells.Item(121, 17).Value = 6
$ws.Cells.Item(121, 18).Value = "Hortaliza"
